# Material link for progress: insert 7 new computer/maths questions
# after the header row, pushing the existing 10 questions down, and
# correct a stray answer reference on the last (English) question.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Make room for the 7 new questions: insert 7 blank rows right
#    after the header (old data rows 2-11 shift down to 9-18).
# ---------------------------------------------------------------
$ws.Rows("2:8").Insert()

# ---------------------------------------------------------------
# 2. Populate the 7 new question rows.
# ---------------------------------------------------------------

# Row 2 - Plug and Play
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Devices that works once they are connected to computer without necessary installing their drivers are called_____"
$ws.Cells.Item(2, 3).Value = "Plug and Play"
$ws.Cells.Item(2, 4).Value = "Plug to Play"
$ws.Cells.Item(2, 5).Value = "Plug then Play"
$ws.Cells.Item(2, 6).Value = "Plug on Play"
$ws.Cells.Item(2, 7).Value = "option_a"

# Row 3 - shortcut key (with rich-text bold on "align center")
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "The short cut keys for align center is __________"
$q3 = $ws.Cells.Item(3, 2)
# touch whole-cell bold once so the bold font gets registered in styles.xml
$q3.Font.Bold = $true
$q3.Font.Bold = $false
$bold = $q3.Characters(24, 12)
$bold.Font.Bold = $true
$bold.Font.Name = "Calibri"
$bold.Font.Size = 11
$tail = $q3.Characters(36, 14)
$tail.Font.Name = "Calibri"
$tail.Font.Size = 11
$ws.Cells.Item(3, 3).Value = "ctlr + H"
$ws.Cells.Item(3, 4).Value = "ctlr + D"
$ws.Cells.Item(3, 5).Value = "ctlr + E"
$ws.Cells.Item(3, 6).Value = "ctlr + AE"
$ws.Cells.Item(3, 7).Value = "option_c"

# Row 4 - application package
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "which of the following is not an example of application package"
$ws.Cells.Item(4, 3).Value = "DYS-E"
$ws.Cells.Item(4, 4).Value = "DBMS"
$ws.Cells.Item(4, 5).Value = "MS-DOS"
$ws.Cells.Item(4, 6).Value = "spread sheet"
$ws.Cells.Item(4, 7).Value = "option_a"

# Row 5 - egyptian hieroglyphics
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "the egyptian where the people that developed _____as symbols representing words"
$ws.Cells.Item(5, 3).Value = "Hieroglyphic"
$ws.Cells.Item(5, 4).Value = "alphabeth"
$ws.Cells.Item(5, 5).Value = "haeroglyphics"
$ws.Cells.Item(5, 6).Value = "algebra"
$ws.Cells.Item(5, 7).Value = "option_a"

# Row 6 - greatest common divisor (numeric options)
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "the greatest common divisor of 126 and 54 is ______"
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 4).Value = 9
$ws.Cells.Item(6, 5).Value = 18
$ws.Cells.Item(6, 6).Value = 54
$ws.Cells.Item(6, 7).Value = "option_c"

# Row 7 - arithmetic / geometric mean
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "find two numbers such that the arithematic mean is 6.5 and geometric mean is 6"
$ws.Cells.Item(7, 3).Value = "3, 6"
$ws.Cells.Item(7, 4).Value = "4, 9"
$ws.Cells.Item(7, 5).Value = "6,8"
$ws.Cells.Item(7, 6).Value = "4, 16"
$ws.Cells.Item(7, 7).Value = "option_b"

# Row 8 - arithmetic series sum (numeric options)
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "if the first and last terms of series are 3 and 35 respectively, find the sum of the first 6 terms"
$ws.Cells.Item(8, 3).Value = 38
$ws.Cells.Item(8, 4).Value = 76
$ws.Cells.Item(8, 5).Value = 114
$ws.Cells.Item(8, 6).Value = 228
$ws.Cells.Item(8, 7).Value = "option_c"

# ---------------------------------------------------------------
# 3. Renumber the "sn" column for the original questions, now
#    shifted down to rows 9-18.
# ---------------------------------------------------------------
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(9 + $i, 1).Value = 8 + $i
}

# ---------------------------------------------------------------
# 4. Fix the stray answer text on the last (English) question -
#    it used to point at a one-off "opton_b" typo string; it
#    should reference the normal "option_b" text instead.
# ---------------------------------------------------------------
$ws.Cells.Item(18, 7).Value = "option_b"

# ---------------------------------------------------------------
# 5. Restore the originally selected cell.
# ---------------------------------------------------------------
$ws.Range("G23").Select()
